$wb = $excel.ActiveWorkbook

# --- Astronauta (sheet 1): new attendance marks in column F (S5) ---
$wsAstro = $wb.Worksheets.Item("Astronauta")
$wsAstro.Range("F5").Value = 1
$wsAstro.Range("F9").Value = 1
$wsAstro.Range("F11").Value = 1
$wsAstro.Range("F23").Value = 1
$wsAstro.Range("F25").Value = 1

# --- Senador (sheet 2): new attendance marks in column F (S5) ---
$wsSenador = $wb.Worksheets.Item("Senador")
$wsSenador.Range("F17").Value = 0
$wsSenador.Range("F23").Value = 0
$wsSenador.Range("F28").Value = 0

# --- Ninja (sheet 4): new attendance marks in column I (T8) ---
$wsNinja = $wb.Worksheets.Item("Ninja")

# Row 2 was recorded as the literal text "0.75" (not a number), so the
# COUNT/SUM formula in column P keeps ignoring it. Build it as a formula
# that yields the text, then flatten to a plain value via copy/paste so no
# new number format / style gets attached to the cell.
$wsNinja.Range("I2").Formula = "=""0.75"""
$wsNinja.Range("I2").Copy()
$wsNinja.Range("I2").PasteSpecial(-4163)

$wsNinja.Range("I5").Value = 1
$wsNinja.Range("I8").Value = 1
$wsNinja.Range("I9").Value = 1
$wsNinja.Range("I11").Value = 1
$wsNinja.Range("I12").Value = 0
$wsNinja.Range("I13").Value = 0
$wsNinja.Range("I14").Value = 1
$wsNinja.Range("I15").Value = 0
$wsNinja.Range("I16").Value = 0
$wsNinja.Range("I17").Value = 1
$wsNinja.Range("I18").Value = 1
$wsNinja.Range("I20").Value = 1
$wsNinja.Range("I22").Value = 1
$wsNinja.Range("I23").Value = 1
$wsNinja.Range("I24").Value = 0
$wsNinja.Range("I25").Value = 1
$wsNinja.Range("I26").Value = 1
$wsNinja.Range("I27").Value = 0
$wsNinja.Range("I28").Value = 0
$wsNinja.Range("I29").Value = 0

# --- Selections / active sheet, matching where each sheet was last left ---
$wsAstro.Range("F25").Select()
$wsNinja.Range("I24").Select()

$wsMago = $wb.Worksheets.Item("Mago")
$wsMago.Range("A24").Select()

$wsSenador.Range("F27").Select()
$wsSenador.Activate()
